$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.400677
$ws.Range("H2").Value = 4.202031
$ws.Range("I2").Value = 0.009106128952548741
$ws.Range("J2").Value = 0.009106128952548741
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03927866666666666
$ws.Range("N2").Value = 0.117836
$ws.Range("O2").Value = 0.7432525340448212
$ws.Range("P2").Value = 0.7432525340448213
$ws.Range("Q2").Value = 0.05501672499066666
$ws.Range("R2").Value = 0.495150524916
$ws.Range("S2").Value = 0.006768153419320765
$ws.Range("T2").Value = 0.006768153419320766

$ws.Range("G3").Value = 1.400677
$ws.Range("H3").Value = 4.202031
$ws.Range("I3").Value = 0.009106128952548741
$ws.Range("J3").Value = 0.009106128952548741
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01356833333333333
$ws.Range("N3").Value = 0.040705
$ws.Range("O3").Value = 0.2567474659551788
$ws.Range("P3").Value = 0.2567474659551788
$ws.Range("Q3").Value = 0.01900485242833333
$ws.Range("R3").Value = 0.171043671855
$ws.Range("S3").Value = 0.002337975533227976
$ws.Range("T3").Value = 0.002337975533227976

$ws.Range("G4").Value = 6.974902666666668
$ws.Range("H4").Value = 20.924708
$ws.Range("I4").Value = 0.04534547444852937
$ws.Range("J4").Value = 0.04534547444852936
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03927866666666666
$ws.Range("N4").Value = 0.117836
$ws.Range("O4").Value = 0.7432525340448212
$ws.Range("P4").Value = 0.7432525340448213
$ws.Range("Q4").Value = 0.2739648768764444
$ws.Range("R4").Value = 2.465683891888
$ws.Range("S4").Value = 0.03370313879133414
$ws.Range("T4").Value = 0.03370313879133414

$ws.Range("G5").Value = 6.974902666666668
$ws.Range("H5").Value = 20.924708
$ws.Range("I5").Value = 0.04534547444852937
$ws.Range("J5").Value = 0.04534547444852936
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01356833333333333
$ws.Range("N5").Value = 0.040705
$ws.Range("O5").Value = 0.2567474659551788
$ws.Range("P5").Value = 0.2567474659551788
$ws.Range("Q5").Value = 0.09463780434888891
$ws.Range("R5").Value = 0.8517402391400001
$ws.Range("S5").Value = 0.01164233565719522
$ws.Range("T5").Value = 0.01164233565719522

$ws.Range("G6").Value = 1.381819666666667
$ws.Range("H6").Value = 4.145459
$ws.Range("I6").Value = 0.00898353301570211
$ws.Range("J6").Value = 0.00898353301570211
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03927866666666666
$ws.Range("N6").Value = 0.117836
$ws.Range("O6").Value = 0.7432525340448212
$ws.Range("P6").Value = 0.7432525340448213
$ws.Range("Q6").Value = 0.05427603408044444
$ws.Range("R6").Value = 0.488484306724
$ws.Range("S6").Value = 0.006677033678595908
$ws.Range("T6").Value = 0.006677033678595908

$ws.Range("G7").Value = 1.381819666666667
$ws.Range("H7").Value = 4.145459
$ws.Range("I7").Value = 0.00898353301570211
$ws.Range("J7").Value = 0.00898353301570211
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01356833333333333
$ws.Range("N7").Value = 0.040705
$ws.Range("O7").Value = 0.2567474659551788
$ws.Range("P7").Value = 0.2567474659551788
$ws.Range("Q7").Value = 0.01874898984388889
$ws.Range("R7").Value = 0.168740908595
$ws.Range("S7").Value = 0.002306499337106202
$ws.Range("T7").Value = 0.002306499337106202

$ws.Range("G8").Value = 144.0595526666667
$ws.Range("H8").Value = 432.178658
$ws.Range("I8").Value = 0.9365648635832199
$ws.Range("J8").Value = 0.9365648635832198
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.03927866666666666
$ws.Range("N8").Value = 0.117836
$ws.Range("O8").Value = 0.7432525340448212
$ws.Range("P8").Value = 0.7432525340448213
$ws.Range("Q8").Value = 5.65846714934311
$ws.Range("R8").Value = 50.926204344088
$ws.Range("S8").Value = 0.6961042081555704
$ws.Range("T8").Value = 0.6961042081555705

$ws.Range("G9").Value = 144.0595526666667
$ws.Range("H9").Value = 432.178658
$ws.Range("I9").Value = 0.9365648635832199
$ws.Range("J9").Value = 0.9365648635832198
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01356833333333333
$ws.Range("N9").Value = 0.040705
$ws.Range("O9").Value = 0.2567474659551788
$ws.Range("P9").Value = 0.2567474659551788
$ws.Range("Q9").Value = 1.954648030432222
$ws.Range("R9").Value = 17.59183227389
$ws.Range("S9").Value = 0.2404606554276494
$ws.Range("T9").Value = 0.2404606554276494

